# New crime data collected — refresh the weekly CompStat figures for the
# 106th Precinct report: bump the report header (volume number + date
# range), then update the Crime Complaints grid (rows 16-29) with the
# newly collected weekly / 28-day / year-to-date / 2-year figures and
# their derived percent-change columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: volume/number and reporting week -----------------------
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# ---- Row 16: Robbery --------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 14
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 209
$ws.Range("J16").Value = 179
$ws.Range("K16").Value = 16.759776536312
$ws.Range("L16").Value = 25.903614457831
$ws.Range("M16").Value = -14.344262295082
$ws.Range("N16").Value = -74.574209245742

# ---- Row 17: Fel. Assault ---------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -45.161290322580
$ws.Range("I17").Value = 320
$ws.Range("J17").Value = 323
$ws.Range("K17").Value = -0.928792569659
$ws.Range("L17").Value = 20.300751879699
$ws.Range("M17").Value = 128.571428571429
$ws.Range("N17").Value = -17.098445595854

# ---- Row 18: Burglary (C18 flips from "n/a" text to a real count) -----
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -72.222222222222
$ws.Range("I18").Value = 138
$ws.Range("J18").Value = 105
$ws.Range("K18").Value = 31.428571428571
$ws.Range("L18").Value = -9.210526315789
$ws.Range("M18").Value = -50.714285714285
$ws.Range("N18").Value = -87.711487088156

# ---- Row 19: Gr. Larceny ----------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 4.166666666666
$ws.Range("I19").Value = 611
$ws.Range("J19").Value = 477
$ws.Range("K19").Value = 28.092243186582
$ws.Range("L19").Value = 22.444889779559
$ws.Range("M19").Value = 80.23598820059
$ws.Range("N19").Value = 11.904761904761

# ---- Row 20: G.L.A. -----------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 207
$ws.Range("J20").Value = 195
$ws.Range("K20").Value = 6.153846153846
$ws.Range("L20").Value = 18.965517241379
$ws.Range("M20").Value = -28.125
$ws.Range("N20").Value = -93.241919686581

# ---- Row 21: TOTAL ------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -21.212121212121
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = -29.861111111111
$ws.Range("I21").Value = 1508
$ws.Range("J21").Value = 1296
$ws.Range("K21").Value = 16.358024691358
$ws.Range("L21").Value = 18.833727344365
$ws.Range("M21").Value = 14.851485148514
$ws.Range("N21").Value = -74.833110814419

# ---- Row 22: Transit (G22/H22 now report "n/a" instead of a value) ----
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("G22").NumberFormat = "General"
$ws.Range("H22").NumberFormat = "General"
$ws.Range("H22").Value = "***.*"

# ---- Row 24: Petit Larceny ---------------------------------------------
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -9.615384615384
$ws.Range("I24").Value = 1264
$ws.Range("J24").Value = 908
$ws.Range("K24").Value = 39.207048458149
$ws.Range("L24").Value = 58.594730238394
$ws.Range("M24").Value = 94.162826420890

# ---- Row 25: Misd. Assault ----------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 497
$ws.Range("J25").Value = 431
$ws.Range("K25").Value = 15.313225058004
$ws.Range("L25").Value = 19.471153846153
$ws.Range("M25").Value = 7.343412526997

# ---- Row 27: Other Sex Crimes -------------------------------------------
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 38
$ws.Range("K27").Value = 39.473684210526
$ws.Range("L27").Value = 82.758620689655

# ---- Row 28: Shooting Vic. (G28/H28 now report "n/a" instead of a value) -
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("H28").Value = "***.*"
$ws.Range("M28").Value = 128.571428571429

# ---- Row 29: Shooting Inc. (G29/H29 now report "n/a" instead of a value) -
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("G29").NumberFormat = "General"
$ws.Range("H29").Value = "***.*"
$ws.Range("M29").Value = 100
